# Insert a new data row before the current row 122, shifting the existing
# rows 122-200 down to 123-201 (dimension grows from A1:R200 to A1:R201),
# then populate the newly inserted row 122 with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above row 122; existing rows 122..200 move to 123..201
$ws.Rows.Item(122).Insert()

# Fill in the brand-new row 122 with the weekly entry that was added.
$ws.Cells.Item(122, 1).Value = 3
$ws.Cells.Item(122, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(122, 3).Value = "Coquimbo"
$ws.Cells.Item(122, 4).Value = 44824
$ws.Cells.Item(122, 5).Value = 5
$ws.Cells.Item(122, 6).Value = 100112026
$ws.Cells.Item(122, 7).Value = "Haba"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 50
$ws.Cells.Item(122, 11).Value = 10000
$ws.Cells.Item(122, 12).Value = 10000
$ws.Cells.Item(122, 13).Value = 10000
$ws.Cells.Item(122, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(122, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(122, 16).Value = 400
$ws.Cells.Item(122, 17).Value = 25
$ws.Cells.Item(122, 18).Value = "Hortaliza"
